# fix: only complain about ambiguous headers if they're used
#
# Adds a pair of ambiguous/duplicate "Duplicate" headers (columns C & D)
# to the products sheet, each with a small 1/2/3 row-number column beneath
# them, so the fixture exercises "ambiguous header present but unused".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): C1/D1 get the new "Duplicate" header -------------
# Clone the formatting of the existing header cell (A1 - bold 13pt font from
# the title-row style) so the new header cells line up with B1/A1, then
# switch the number format to General (headers aren't forced to text here,
# matching the style already used internally by the workbook: xf index 5).
$ws.Range("A1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$ws.Range("C1").NumberFormat = "General"
$ws.Range("C1").Value = "Duplicate"

$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("D1").Value = "Duplicate"

# --- Data rows: sequential numbers 1, 2, 3 under each "Duplicate" column --
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 1
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 2
$ws.Range("C4").Value = 3
$ws.Range("D4").Value = 3

# --- Move the selection cursor to D2 (matches the authored selection) -----
$ws.Range("D2").Select() | Out-Null
